$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.1926910299003322
$ws.Range("C2").Value = 0.5681063122923588
$ws.Range("J2").Value = 0.02325581395348837
$ws.Range("P2").Value = 0.1395348837209302
$ws.Range("S2").Value = 0.07641196013289037
# Row 3
$ws.Range("B3").Value = 0.005780346820809248
$ws.Range("C3").Value = 0.01734104046242774
$ws.Range("J3").Value = 0.02312138728323699
$ws.Range("P3").Value = 0.7630057803468208
$ws.Range("S3").Value = 0.1907514450867052
# Row 4
$ws.Range("J4").Value = 0.03278688524590164
$ws.Range("P4").Value = 0.7704918032786885
$ws.Range("S4").Value = 0.1967213114754098
# Row 6
$ws.Range("B6").Value = 0.07391304347826087
$ws.Range("D6").Value = 0.01304347826086956
$ws.Range("F6").Value = 0.06956521739130435
$ws.Range("J6").Value = 0.2043478260869565
$ws.Range("O6").Value = 0.004347826086956522
$ws.Range("Q6").Value = 0.2
$ws.Range("R6").Value = 0.0782608695652174
$ws.Range("S6").Value = 0.3565217391304348
# Row 7
$ws.Range("B7").Value = 0.1125
$ws.Range("D7").Value = 0.01875
$ws.Range("F7").Value = 0.05625
$ws.Range("J7").Value = 0.15625
$ws.Range("O7").Value = 0.03125
$ws.Range("Q7").Value = 0.23125
$ws.Range("R7").Value = 0.0375
$ws.Range("S7").Value = 0.35625
# Row 8
$ws.Range("B8").Value = 0.0995260663507109
$ws.Range("D8").Value = 0.01895734597156398
$ws.Range("F8").Value = 0.07109004739336493
$ws.Range("J8").Value = 0.1303317535545024
$ws.Range("O8").Value = 0.01421800947867299
$ws.Range("Q8").Value = 0.2085308056872038
$ws.Range("R8").Value = 0.06872037914691943
$ws.Range("S8").Value = 0.3886255924170616
# Row 9
$ws.Range("B9").Value = 0.1117647058823529
$ws.Range("D9").Value = 0.04705882352941176
$ws.Range("E9").Value = 0.005882352941176471
$ws.Range("F9").Value = 0.09411764705882353
$ws.Range("J9").Value = 0.09411764705882353
$ws.Range("O9").Value = 0.02941176470588235
$ws.Range("Q9").Value = 0.1529411764705882
$ws.Range("R9").Value = 0.09411764705882353
$ws.Range("S9").Value = 0.3705882352941177
# Row 10
$ws.Range("B10").Value = 0.1094946401225115
$ws.Range("D10").Value = 0.03215926493108729
$ws.Range("E10").Value = 0.002297090352220521
$ws.Range("F10").Value = 0.07044410413476264
$ws.Range("J10").Value = 0.1309341500765697
$ws.Range("O10").Value = 0.01761102603369066
$ws.Range("Q10").Value = 0.217457886676876
$ws.Range("R10").Value = 0.08652373660030628
$ws.Range("S10").Value = 0.3330781010719755
# Row 11
$ws.Range("F11").Value = 0.003891050583657588
$ws.Range("G11").Value = 0.1478599221789883
$ws.Range("J11").Value = 0.1011673151750973
$ws.Range("K11").Value = 0.2217898832684825
$ws.Range("L11").Value = 0.5058365758754864
$ws.Range("S11").Value = 0.01945525291828794
# Row 12
$ws.Range("G12").Value = 0.7313432835820896
$ws.Range("J12").Value = 0.1865671641791045
$ws.Range("K12").Value = 0.02238805970149254
$ws.Range("L12").Value = 0.02238805970149254
$ws.Range("S12").Value = 0.03731343283582089
# Row 13
$ws.Range("G13").Value = 0.5957446808510638
$ws.Range("J13").Value = 0.3191489361702128
$ws.Range("S13").Value = 0.0851063829787234
# Row 14
$ws.Range("G14").Value = 0.5
$ws.Range("J14").Value = 0.5
# Row 15
$ws.Range("F15").Value = 0.01481481481481482
$ws.Range("H15").Value = 0.1333333333333333
$ws.Range("I15").Value = 0.07037037037037037
$ws.Range("J15").Value = 0.3962962962962963
$ws.Range("K15").Value = 0.05555555555555555
$ws.Range("M15").Value = 0.007407407407407408
$ws.Range("N15").Value = 0.003703703703703704
$ws.Range("O15").Value = 0.1111111111111111
$ws.Range("S15").Value = 0.2074074074074074
# Row 16
$ws.Range("F16").Value = 0.009302325581395349
$ws.Range("H16").Value = 0.1674418604651163
$ws.Range("I16").Value = 0.08372093023255814
$ws.Range("J16").Value = 0.4
$ws.Range("K16").Value = 0.1162790697674419
$ws.Range("M16").Value = 0.004651162790697674
$ws.Range("O16").Value = 0.08837209302325581
$ws.Range("S16").Value = 0.1302325581395349
# Row 17
$ws.Range("F17").Value = 0.006355932203389831
$ws.Range("H17").Value = 0.173728813559322
$ws.Range("I17").Value = 0.06991525423728813
$ws.Range("J17").Value = 0.4661016949152542
$ws.Range("K17").Value = 0.0826271186440678
$ws.Range("M17").Value = 0.01059322033898305
$ws.Range("O17").Value = 0.06567796610169492
$ws.Range("S17").Value = 0.125
# Row 18
$ws.Range("F18").Value = 0.01666666666666667
$ws.Range("I18").Value = 0.07777777777777778
$ws.Range("J18").Value = 0.4722222222222222
$ws.Range("K18").Value = 0.08888888888888889
$ws.Range("M18").Value = 0.005555555555555556
$ws.Range("N18").Value = 0.005555555555555556
$ws.Range("S18").Value = 0.1333333333333333
# Row 19
$ws.Range("F19").Value = 0.01957446808510638
$ws.Range("H19").Value = 0.2059574468085106
$ws.Range("I19").Value = 0.0748936170212766
$ws.Range("J19").Value = 0.3676595744680851
$ws.Range("K19").Value = 0.08425531914893616
$ws.Range("M19").Value = 0.03404255319148936
$ws.Range("O19").Value = 0.0902127659574468
$ws.Range("S19").Value = 0.1234042553191489
